$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.370.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.66%  "

$ws.Range("D3").Value = "'1.940.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.61%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'250.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.29%  "

$ws.Range("D6").Value = "'0.7259"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.32%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").Value = "'0.3338"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.74%  "

$ws.Range("D9").Value = "'28.62"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.91%  "

$ws.Range("D10").Value = "'0.07264"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.35%  "

$ws.Range("D11").Value = "'0.8130"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.57%  "

$ws.Range("D12").Value = "'0.08105"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.88%  "

$ws.Range("D13").Value = "'1.937.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.77%  "

$ws.Range("D14").Value = "'5.487"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.36%  "

$ws.Range("D15").Value = "'94.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.66%  "

$ws.Range("D16").Value = "'15.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.18%  "

$ws.Range("D17").Value = "'30.354.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.75%  "

$ws.Range("D18").Value = "'0.000008219"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.56%  "

$ws.Range("D19").Value = "'252.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.35%  "

$ws.Range("D20").Value = "'5.931"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.14%  "

$ws.Range("D21").Value = "'2.189.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.93%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").Value = "'1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").Value = "'6.960"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.22%  "

$ws.Range("D25").Value = "'9.787"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.81%  "

$ws.Range("D26").Value = "'163.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.09%  "

$ws.Range("D27").Value = "'2.401"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.38%  "

$ws.Range("D28").Value = "'19.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.38%  "

$ws.Range("D29").Value = "'0.1324"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.42%  "

$ws.Range("D30").Value = "'1.573"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.24%  "

$ws.Range("D31").Value = "'1.347"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.09%  "

$ws.Range("D32").Value = "'4.452"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.94%  "

$ws.Range("D33").Value = "'4.215"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.76%  "

$ws.Range("D34").Value = "'0.05201"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.73%  "

$ws.Range("D35").Value = "'1.299"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.56%  "

$ws.Range("D36").Value = "'0.7527"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.93%  "

$ws.Range("D37").Value = "'2.749"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.64%  "

$ws.Range("D38").Value = "'0.01982"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.06%  "

$ws.Range("D39").Value = "'2.838"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.00%  "

$ws.Range("D40").Value = "'81.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.13%  "

$ws.Range("D41").Value = "'6.541"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.30%  "

$ws.Range("D42").Value = "'0.4552"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.51%  "

$ws.Range("D43").Value = "'2.046"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.33%  "

$ws.Range("D44").Value = "'0.8481"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.46%  "

$ws.Range("D45").Value = "'1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("D46").Value = "'102.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.98%  "

$ws.Range("D47").Value = "'9.803"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.71%  "

$ws.Range("D48").Value = "'7.468"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.39%  "

$ws.Range("D49").Value = "'36.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.92%  "

$ws.Range("D50").Value = "'0.4202"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.32%  "

$ws.Range("D51").Value = "'0.06055"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.26%  "
